# texte soutenance + orthographe
#
# Re-types a handful of text runs that had been split into several
# a:r runs (e.g. because of autocorrect / incremental typing) back into a
# single run, and fixes a couple of spelling / wording mistakes across the
# "soutenance" slides. Also clears the empty leftover <p:timing> nodes on
# the slides we touch, the same way PowerPoint drops them once you resave
# a slide that has no real animations on it.

$p = $ppt.ActivePresentation

function Clear-EmptyTiming($slide) {
    $seq = $slide.TimeLine.MainSequence
    $effect = $seq.AddEffect($slide.Shapes.Item(1), 1)
    $effect.Delete()
}

function Set-ShapeText($shape, [string]$text) {
    $shape.TextFrame.TextRange.Delete()
    $shape.TextFrame.TextRange.InsertAfter($text) | Out-Null
}

# --- Slide 1 : title slide -------------------------------------------------
$s1 = $p.Slides.Item(1)
Clear-EmptyTiming $s1

# --- Slide 3 : "IDENTIFIER ET DEFINIR LE PROBLEME" --------------------------
$s3 = $p.Slides.Item(3)
$nbsp = [char]0x00A0
Set-ShapeText $s3.Shapes.Item(5) "Pouvons-nous détecter des cellules cancéreuses sur une coupe d’échantillon de tumeur afin de localiser précisément les zones cancéreuses et ainsi évaluer la gravité du cancer$($nbsp)?"
Set-ShapeText $s3.Shapes.Item(7) "Implémentation de Machine Learning"
Clear-EmptyTiming $s3

# --- Slide 4 : "DESCRIPTION DE LA BASE DE DONNEES" --------------------------
$s4 = $p.Slides.Item(4)
Set-ShapeText $s4.Shapes.Item(7) "277 524 images"
Set-ShapeText $s4.Shapes.Item(8) "275 246 images"
Clear-EmptyTiming $s4

# --- Slide 5 : "CHARGER LA BASE DE DONNEES DANS UNE STRUCTURE" (Data frame) -
$s5 = $p.Slides.Item(5)
Clear-EmptyTiming $s5

# --- Slide 6 : "CHARGER LA BASE DE DONNEES DANS UNE STRUCTURE" (Tableau) ----
$s6 = $p.Slides.Item(6)
Clear-EmptyTiming $s6

# --- Slide 7 : "PREPARATION A L'APPLICATION" --------------------------------
$s7 = $p.Slides.Item(7)
Clear-EmptyTiming $s7

# --- Slide 8 : "APLICATION DES MODELES AUX DONNEES" (KNN) -------------------
$s8 = $p.Slides.Item(8)
Set-ShapeText $s8.Shapes.Item(2) "APPLICATION DES MODELES AUX DONNEES"
Set-ShapeText $s8.Shapes.Item(5) "Temps d’exécution = 8 heures"
Clear-EmptyTiming $s8

# --- Slide 9 : "APLICATION DES MODELES AUX DONNEES" (arbres de décision) ----
$s9 = $p.Slides.Item(9)
Set-ShapeText $s9.Shapes.Item(2) "2 – Classification avec arbres de décision"
Set-ShapeText $s9.Shapes.Item(3) "APPLICATION DES MODELES AUX DONNEES"
Set-ShapeText $s9.Shapes.Item(8) "Temps d’exécution = 20 minutes"
Clear-EmptyTiming $s9

# --- Slide 10 : "APLICATION DES MODELES AUX DONNEES" (SVM) ------------------
$s10 = $p.Slides.Item(10)
Set-ShapeText $s10.Shapes.Item(2) "3 – Classification avec SVM"
Set-ShapeText $s10.Shapes.Item(3) "APPLICATION DES MODELES AUX DONNEES"
Set-ShapeText $s10.Shapes.Item(5) "Temps d’exécution = 11 heures"
Clear-EmptyTiming $s10

# --- Slide 11 : "CONCLUSION" -------------------------------------------------
$s11 = $p.Slides.Item(11)
Clear-EmptyTiming $s11
